$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 79; this shifts the existing rows 79-92 down to 80-93,
# matching the dimension change from A1:R92 to A1:R93.
$ws.Rows.Item(79).Insert()

# Populate the newly inserted row 79 with the new weekly price record.
$ws.Cells.Item(79, 1).Value2 = 10
$ws.Cells.Item(79, 2).Value2 = "Vega Modelo de Temuco"
$ws.Cells.Item(79, 3).Value2 = "La Araucanía"
$ws.Cells.Item(79, 4).Value2 = 45204
$ws.Cells.Item(79, 5).Value2 = 9
$ws.Cells.Item(79, 6).Value2 = 300000000
$ws.Cells.Item(79, 7).Value2 = "Espárragos"
$ws.Cells.Item(79, 8).Value2 = "Sin especificar"
$ws.Cells.Item(79, 9).Value2 = "Primera"
$ws.Cells.Item(79, 10).Value2 = 680
$ws.Cells.Item(79, 11).Value2 = 1500
$ws.Cells.Item(79, 12).Value2 = 1600
$ws.Cells.Item(79, 13).Value2 = 1571
$ws.Cells.Item(79, 14).Value2 = "`$/kilo"
$ws.Cells.Item(79, 15).Value2 = "Región del Maule"
$ws.Cells.Item(79, 16).Value2 = 1571
$ws.Cells.Item(79, 17).Value2 = 1
$ws.Cells.Item(79, 18).Value2 = "Hortaliza"
